# Applies scheduled-runner market-data refresh to the 8 crafting-leve sheets
# (currentAveragePrice* / LevePrice* / LeveProfit* columns), per the upstream diff.
$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

# ALC!40 (Leve Item ID 5505)
$wsALC.Cells.Item(40, 8).Value = 3370.4614
$wsALC.Cells.Item(40, 9).Value = 2266.6667
$wsALC.Cells.Item(40, 10).Value = 3701.6
$wsALC.Cells.Item(40, 11).Value = 2266.6667
$wsALC.Cells.Item(40, 12).Value = 3701.6
$wsALC.Cells.Item(40, 13).Value = -2091.6667
$wsALC.Cells.Item(40, 14).Value = -4051.6

# ALC!55 (Leve Item ID 5517)
$wsALC.Cells.Item(55, 8).Value = 343.66666
$wsALC.Cells.Item(55, 10).Value = 2
$wsALC.Cells.Item(55, 12).Value = 2
$wsALC.Cells.Item(55, 14).Value = -430

# ALC!64 (Leve Item ID 5506)
$wsALC.Cells.Item(64, 8).Value = 8000
$wsALC.Cells.Item(64, 9).Value = 8000
$wsALC.Cells.Item(64, 11).Value = 8000
$wsALC.Cells.Item(64, 13).Value = -7752

# ALC!67 (Leve Item ID 5506)
$wsALC.Cells.Item(67, 8).Value = 8000
$wsALC.Cells.Item(67, 9).Value = 8000
$wsALC.Cells.Item(67, 11).Value = 8000
$wsALC.Cells.Item(67, 13).Value = -7142

# ALC!74 (Leve Item ID 5507)
$wsALC.Cells.Item(74, 8).Value = 4700
$wsALC.Cells.Item(74, 9).Value = 4700
$wsALC.Cells.Item(74, 11).Value = 4700
$wsALC.Cells.Item(74, 13).Value = -3764

# ALC!77 (Leve Item ID 5507)
$wsALC.Cells.Item(77, 8).Value = 4700
$wsALC.Cells.Item(77, 9).Value = 4700
$wsALC.Cells.Item(77, 11).Value = 23500
$wsALC.Cells.Item(77, 13).Value = -18820

# ARM!63 (Leve Item ID 12528)
$wsARM.Cells.Item(63, 8).Value = 5235.5713
$wsARM.Cells.Item(63, 9).Value = 5729.8
$wsARM.Cells.Item(63, 11).Value = 5729.8
$wsARM.Cells.Item(63, 13).Value = -5043.8

# ARM!66 (Leve Item ID 12528)
$wsARM.Cells.Item(66, 8).Value = 5235.5713
$wsARM.Cells.Item(66, 9).Value = 5729.8
$wsARM.Cells.Item(66, 11).Value = 28649
$wsARM.Cells.Item(66, 13).Value = -25217

# ARM!122 (Leve Item ID 36168)
$wsARM.Cells.Item(122, 8).Value = 3753
$wsARM.Cells.Item(122, 9).Value = 3753
$wsARM.Cells.Item(122, 10).Value = 0
$wsARM.Cells.Item(122, 11).Value = 11259
$wsARM.Cells.Item(122, 12).Value = 0
$wsARM.Cells.Item(122, 13).Value = -8809
$wsARM.Cells.Item(122, 14).Value = ""

# ARM!132 (Leve Item ID 43997)
$wsARM.Cells.Item(132, 8).Value = 4105.4165
$wsARM.Cells.Item(132, 9).Value = 1362.7778
$wsARM.Cells.Item(132, 10).Value = 12333.333
$wsARM.Cells.Item(132, 11).Value = 4088.3334
$wsARM.Cells.Item(132, 12).Value = 36999.999
$wsARM.Cells.Item(132, 13).Value = -1558.3334
$wsARM.Cells.Item(132, 14).Value = -42059.999

# BSM!102 (Leve Item ID 19565)
$wsBSM.Cells.Item(102, 8).Value = 37333.332
$wsBSM.Cells.Item(102, 9).Value = 18500
$wsBSM.Cells.Item(102, 11).Value = 18500
$wsBSM.Cells.Item(102, 13).Value = -15255

# BSM!105 (Leve Item ID 19947)
$wsBSM.Cells.Item(105, 8).Value = 2195.4285
$wsBSM.Cells.Item(105, 9).Value = 2144.6667
$wsBSM.Cells.Item(105, 11).Value = 2144.6667
$wsBSM.Cells.Item(105, 13).Value = -397.6667000000002

# BSM!134 (Leve Item ID 43998)
$wsBSM.Cells.Item(134, 8).Value = 3082.8235
$wsBSM.Cells.Item(134, 9).Value = 1416.0769
$wsBSM.Cells.Item(134, 11).Value = 4248.2307
$wsBSM.Cells.Item(134, 13).Value = -1713.2307

# CRP!17 (Leve Item ID 1823)
$wsCRP.Cells.Item(17, 8).Value = 0
$wsCRP.Cells.Item(17, 9).Value = 0
$wsCRP.Cells.Item(17, 11).Value = 0
$wsCRP.Cells.Item(17, 13).Value = ""

# CRP!22 (Leve Item ID 5367)
$wsCRP.Cells.Item(22, 8).Value = 280
$wsCRP.Cells.Item(22, 9).Value = 280
$wsCRP.Cells.Item(22, 11).Value = 280
$wsCRP.Cells.Item(22, 13).Value = 70

# CRP!62 (Leve Item ID 12580)
$wsCRP.Cells.Item(62, 8).Value = 0
$wsCRP.Cells.Item(62, 9).Value = 0
$wsCRP.Cells.Item(62, 11).Value = 0
$wsCRP.Cells.Item(62, 13).Value = ""

# CRP!65 (Leve Item ID 12580)
$wsCRP.Cells.Item(65, 8).Value = 0
$wsCRP.Cells.Item(65, 9).Value = 0
$wsCRP.Cells.Item(65, 11).Value = 0
$wsCRP.Cells.Item(65, 13).Value = ""

# CRP!141 (Leve Item ID 43345)
$wsCRP.Cells.Item(141, 8).Value = 480101.38
$wsCRP.Cells.Item(141, 10).Value = 480101.38
$wsCRP.Cells.Item(141, 12).Value = 480101.38
$wsCRP.Cells.Item(141, 14).Value = -490461.38

# CUL!5 (Leve Item ID 43974)
$wsCUL.Cells.Item(5, 8).Value = 3543.75
$wsCUL.Cells.Item(5, 9).Value = 485
$wsCUL.Cells.Item(5, 11).Value = 1455
$wsCUL.Cells.Item(5, 13).Value = -1343

# CUL!14 (Leve Item ID 12886)
$wsCUL.Cells.Item(14, 8).Value = 3048.7144
$wsCUL.Cells.Item(14, 9).Value = 3048.7144
$wsCUL.Cells.Item(14, 11).Value = 9146.143199999999
$wsCUL.Cells.Item(14, 13).Value = -8973.143199999999

# CUL!34 (Leve Item ID 4749)
$wsCUL.Cells.Item(34, 8).Value = 787.2222
$wsCUL.Cells.Item(34, 9).Value = 425
$wsCUL.Cells.Item(34, 10).Value = 1077
$wsCUL.Cells.Item(34, 11).Value = 1275
$wsCUL.Cells.Item(34, 12).Value = 3231
$wsCUL.Cells.Item(34, 13).Value = -1191
$wsCUL.Cells.Item(34, 14).Value = -3399

# CUL!39 (Leve Item ID 4712)
$wsCUL.Cells.Item(39, 8).Value = 3667.889
$wsCUL.Cells.Item(39, 10).Value = 3667.889
$wsCUL.Cells.Item(39, 12).Value = 11003.667
$wsCUL.Cells.Item(39, 14).Value = -11591.667

# CUL!135 (Leve Item ID 43974)
$wsCUL.Cells.Item(135, 8).Value = 3543.75
$wsCUL.Cells.Item(135, 9).Value = 485
$wsCUL.Cells.Item(135, 11).Value = 4365
$wsCUL.Cells.Item(135, 13).Value = -1830

# GSM!2 (Leve Item ID 5062)
$wsGSM.Cells.Item(2, 8).Value = 60.46154
$wsGSM.Cells.Item(2, 9).Value = 35.75
$wsGSM.Cells.Item(2, 11).Value = 35.75
$wsGSM.Cells.Item(2, 13).Value = 77.25

# GSM!10 (Leve Item ID 4306)
$wsGSM.Cells.Item(10, 8).Value = 3475.5
$wsGSM.Cells.Item(10, 10).Value = 1299.6666
$wsGSM.Cells.Item(10, 12).Value = 1299.6666
$wsGSM.Cells.Item(10, 14).Value = -1637.6666

# GSM!46 (Leve Item ID 2078)
$wsGSM.Cells.Item(46, 8).Value = 0
$wsGSM.Cells.Item(46, 9).Value = 0
$wsGSM.Cells.Item(46, 11).Value = 0
$wsGSM.Cells.Item(46, 13).Value = ""

# GSM!48 (Leve Item ID 4337)
$wsGSM.Cells.Item(48, 8).Value = 0
$wsGSM.Cells.Item(48, 10).Value = 0
$wsGSM.Cells.Item(48, 12).Value = 0
$wsGSM.Cells.Item(48, 14).Value = ""

# GSM!58 (Leve Item ID 4363)
$wsGSM.Cells.Item(58, 8).Value = 25857
$wsGSM.Cells.Item(58, 10).Value = 31499.75
$wsGSM.Cells.Item(58, 12).Value = 31499.75
$wsGSM.Cells.Item(58, 14).Value = -32053.75

# GSM!80 (Leve Item ID 12521)
$wsGSM.Cells.Item(80, 8).Value = 3999.5
$wsGSM.Cells.Item(80, 9).Value = 0
$wsGSM.Cells.Item(80, 11).Value = 0
$wsGSM.Cells.Item(80, 13).Value = ""

# GSM!83 (Leve Item ID 12521)
$wsGSM.Cells.Item(83, 8).Value = 3999.5
$wsGSM.Cells.Item(83, 9).Value = 0
$wsGSM.Cells.Item(83, 11).Value = 0
$wsGSM.Cells.Item(83, 13).Value = ""

# LTW!16 (Leve Item ID 5289)
$wsLTW.Cells.Item(16, 8).Value = 1650.3334
$wsLTW.Cells.Item(16, 9).Value = 1650.3334
$wsLTW.Cells.Item(16, 11).Value = 1650.3334
$wsLTW.Cells.Item(16, 13).Value = -1480.3334

# LTW!46 (Leve Item ID 5282)
$wsLTW.Cells.Item(46, 8).Value = 6564.3125
$wsLTW.Cells.Item(46, 10).Value = 5016.25
$wsLTW.Cells.Item(46, 12).Value = 5016.25
$wsLTW.Cells.Item(46, 14).Value = -5392.25

# LTW!68 (Leve Item ID 12563)
$wsLTW.Cells.Item(68, 8).Value = 3193.3333
$wsLTW.Cells.Item(68, 9).Value = 3192
$wsLTW.Cells.Item(68, 11).Value = 3192
$wsLTW.Cells.Item(68, 13).Value = -2443

# LTW!71 (Leve Item ID 12563)
$wsLTW.Cells.Item(71, 8).Value = 3193.3333
$wsLTW.Cells.Item(71, 9).Value = 3192
$wsLTW.Cells.Item(71, 11).Value = 15960
$wsLTW.Cells.Item(71, 13).Value = -12216

# LTW!82 (Leve Item ID 12565)
$wsLTW.Cells.Item(82, 8).Value = 1300.8889
$wsLTW.Cells.Item(82, 9).Value = 902.6667
$wsLTW.Cells.Item(82, 10).Value = 1500
$wsLTW.Cells.Item(82, 11).Value = 902.6667
$wsLTW.Cells.Item(82, 12).Value = 1500
$wsLTW.Cells.Item(82, 13).Value = -541.6667
$wsLTW.Cells.Item(82, 14).Value = -2222

# LTW!85 (Leve Item ID 12565)
$wsLTW.Cells.Item(85, 8).Value = 1300.8889
$wsLTW.Cells.Item(85, 9).Value = 902.6667
$wsLTW.Cells.Item(85, 10).Value = 1500
$wsLTW.Cells.Item(85, 11).Value = 902.6667
$wsLTW.Cells.Item(85, 12).Value = 1500
$wsLTW.Cells.Item(85, 13).Value = 345.3333
$wsLTW.Cells.Item(85, 14).Value = -3996

# LTW!93 (Leve Item ID 19993)
$wsLTW.Cells.Item(93, 8).Value = 1896.7142
$wsLTW.Cells.Item(93, 9).Value = 1655.4
$wsLTW.Cells.Item(93, 11).Value = 1655.4
$wsLTW.Cells.Item(93, 13).Value = -407.4000000000001

# WVR!4 (Leve Item ID 2996)
$wsWVR.Cells.Item(4, 8).Value = 34900
$wsWVR.Cells.Item(4, 9).Value = 50000
$wsWVR.Cells.Item(4, 10).Value = 19800
$wsWVR.Cells.Item(4, 11).Value = 50000
$wsWVR.Cells.Item(4, 12).Value = 19800
$wsWVR.Cells.Item(4, 13).Value = -49887
$wsWVR.Cells.Item(4, 14).Value = -20026

# WVR!126 (Leve Item ID 36210)
$wsWVR.Cells.Item(126, 8).Value = 0
$wsWVR.Cells.Item(126, 9).Value = 0
$wsWVR.Cells.Item(126, 10).Value = 0
$wsWVR.Cells.Item(126, 11).Value = 0
$wsWVR.Cells.Item(126, 12).Value = 0
$wsWVR.Cells.Item(126, 13).Value = ""
$wsWVR.Cells.Item(126, 14).Value = ""

# WVR!132 (Leve Item ID 44029)
$wsWVR.Cells.Item(132, 8).Value = 6599.4
$wsWVR.Cells.Item(132, 9).Value = 999
$wsWVR.Cells.Item(132, 10).Value = 15000
$wsWVR.Cells.Item(132, 11).Value = 2997
$wsWVR.Cells.Item(132, 12).Value = 45000
$wsWVR.Cells.Item(132, 13).Value = -467
$wsWVR.Cells.Item(132, 14).Value = -50060
